$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 174.76923
$ws.Range("I38").Value = 174.76923
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 524.30769
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -152.30769
$ws.Range("N38").ClearContents()

$ws.Range("H40").Value = 6515.2856
$ws.Range("J40").Value = 9799.200000000001
$ws.Range("L40").Value = 9799.200000000001
$ws.Range("N40").Value = -10149.2

$ws.Range("H41").Value = 978.3333
$ws.Range("I41").Value = 174
$ws.Range("K41").Value = 174
$ws.Range("M41").Value = 266

$ws.Range("H62").Value = 10500.125
$ws.Range("I62").Value = 10000.167
$ws.Range("K62").Value = 10000.167
$ws.Range("M62").Value = -9376.166999999999

$ws.Range("H65").Value = 10500.125
$ws.Range("I65").Value = 10000.167
$ws.Range("K65").Value = 50000.835
$ws.Range("M65").Value = -46880.835

$ws.Range("H98").Value = 2087.9
$ws.Range("I98").Value = 1764.3334
$ws.Range("K98").Value = 1764.3334
$ws.Range("M98").Value = -266.3334

$ws.Range("H122").Value = 2087.9
$ws.Range("I122").Value = 1764.3334
$ws.Range("K122").Value = 5293.0002
$ws.Range("M122").Value = -2843.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 497.33334
$ws.Range("I26").Value = 497.33334
$ws.Range("K26").Value = 497.33334
$ws.Range("M26").Value = -167.33334

$ws.Range("H28").Value = 14842.5
$ws.Range("I28").Value = 14842.5
$ws.Range("K28").Value = 14842.5
$ws.Range("M28").Value = -14650.5

$ws.Range("H38").Value = 14954.2
$ws.Range("I38").Value = 2937.5
$ws.Range("K38").Value = 2937.5
$ws.Range("M38").Value = -2470.5

$ws.Range("H39").Value = 6591.8823
$ws.Range("I39").Value = 3065.25
$ws.Range("J39").Value = 63018
$ws.Range("K39").Value = 3065.25
$ws.Range("L39").Value = 63018
$ws.Range("M39").Value = -2545.25
$ws.Range("N39").Value = -64058

$ws.Range("H74").Value = 5564.8
$ws.Range("I74").Value = 4466.0835
$ws.Range("K74").Value = 4466.0835
$ws.Range("M74").Value = -3592.0835

$ws.Range("H77").Value = 5564.8
$ws.Range("I77").Value = 4466.0835
$ws.Range("K77").Value = 22330.4175
$ws.Range("M77").Value = -17962.4175

$ws.Range("H99").Value = 14842.5
$ws.Range("I99").Value = 14842.5
$ws.Range("K99").Value = 14842.5
$ws.Range("M99").Value = -11847.5

$ws.Range("H122").Value = 1806.1428
$ws.Range("I122").Value = 1578.6
$ws.Range("K122").Value = 4735.799999999999
$ws.Range("M122").Value = -2285.799999999999

$ws.Range("H137").Value = 74999
$ws.Range("J137").Value = 74999
$ws.Range("L137").Value = 74999
$ws.Range("N137").Value = -85199

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 18623
$ws.Range("J88").Value = 18623
$ws.Range("L88").Value = 18623
$ws.Range("N88").Value = -19435

$ws.Range("H91").Value = 18623
$ws.Range("J91").Value = 18623
$ws.Range("L91").Value = 18623
$ws.Range("N91").Value = -21431

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("N130").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7423.25
$ws.Range("I6").Value = 340.57144
$ws.Range("J6").Value = 57002
$ws.Range("K6").Value = 340.57144
$ws.Range("L6").Value = 57002
$ws.Range("M6").Value = -227.57144
$ws.Range("N6").Value = -57228

$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H18").Value = 49990
$ws.Range("J18").Value = 49990
$ws.Range("L18").Value = 49990
$ws.Range("N18").Value = -50450

$ws.Range("H32").Value = 1130
$ws.Range("I32").Value = 1130
$ws.Range("K32").Value = 1130
$ws.Range("M32").Value = -814

$ws.Range("H38").Value = 12521
$ws.Range("I38").Value = 5000
$ws.Range("J38").Value = 20042
$ws.Range("K38").Value = 5000
$ws.Range("L38").Value = 20042
$ws.Range("M38").Value = -4623
$ws.Range("N38").Value = -20796

$ws.Range("H46").Value = 12521
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 20042
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 20042
$ws.Range("M46").Value = -4789
$ws.Range("N46").Value = -20464

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 1438
$ws.Range("I122").Value = 1290.3889
$ws.Range("K122").Value = 3871.1667
$ws.Range("M122").Value = -1421.1667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 16499.5
$ws.Range("J88").Value = 16499.5
$ws.Range("L88").Value = 49498.5
$ws.Range("N88").Value = -50354.5

$ws.Range("H91").Value = 16499.5
$ws.Range("J91").Value = 16499.5
$ws.Range("L91").Value = 49498.5
$ws.Range("N91").Value = -52462.5

$ws.Range("H128").Value = 1849999.2
$ws.Range("I128").Value = 1849999.2
$ws.Range("K128").Value = 5549997.6
$ws.Range("M128").Value = -5545017.6

$ws.Range("H136").Value = 6832.6665
$ws.Range("I136").Value = 499
$ws.Range("K136").Value = 1497
$ws.Range("M136").Value = 3603

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4725
$ws.Range("I70").Value = 4725
$ws.Range("K70").Value = 4725
$ws.Range("M70").Value = -4455

$ws.Range("H73").Value = 4725
$ws.Range("I73").Value = 4725
$ws.Range("K73").Value = 4725
$ws.Range("M73").Value = -3789

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H102").Value = 1402.129
$ws.Range("I102").Value = 869.6667
$ws.Range("K102").Value = 869.6667
$ws.Range("M102").Value = 752.3333

$ws.Range("H122").Value = 211662.33
$ws.Range("I122").Value = 297008.2
$ws.Range("J122").Value = 4393.857
$ws.Range("K122").Value = 891024.6000000001
$ws.Range("L122").Value = 13181.571
$ws.Range("M122").Value = -888574.6000000001
$ws.Range("N122").Value = -18081.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 5000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -5590

$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 5000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -5214

$ws.Range("H40").Value = 3926.4075
$ws.Range("I40").Value = 3001.1738
$ws.Range("J40").Value = 9246.5
$ws.Range("K40").Value = 3001.1738
$ws.Range("L40").Value = 9246.5
$ws.Range("M40").Value = -2865.1738
$ws.Range("N40").Value = -9518.5

$ws.Range("H98").Value = 51077.668
$ws.Range("J98").Value = 51077.668
$ws.Range("L98").Value = 51077.668
$ws.Range("N98").Value = -57067.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 34866.5
$ws.Range("J98").Value = 34866.5
$ws.Range("L98").Value = 34866.5
$ws.Range("N98").Value = -40856.5

$ws.Range("H107").Value = 488.57144
$ws.Range("J107").Value = 397
$ws.Range("L107").Value = 1191
$ws.Range("N107").Value = -5031

$ws.Range("H122").Value = 1580.7273
$ws.Range("I122").Value = 1580.7273
$ws.Range("K122").Value = 4742.1819
$ws.Range("M122").Value = -2292.1819

$ws.Range("H126").Value = 3024.2
$ws.Range("I126").Value = 1496.7646
$ws.Range("K126").Value = 4490.293799999999
$ws.Range("M126").Value = -2020.293799999999

$ws.Range("H127").Value = 222500
$ws.Range("I127").Value = 220000
$ws.Range("K127").Value = 220000
$ws.Range("M127").Value = -215040

$ws.Range("H136").Value = 2678.6943
$ws.Range("I136").Value = 1609.5172
$ws.Range("K136").Value = 4828.5516
$ws.Range("M136").Value = -2678.5516
